$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to Text format so numeric-looking strings
# (e.g. trailing-zero decimals, multi-dot thousand separators) are preserved
# exactly as text, matching the source inlineStr cells.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.535.77'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.933.88'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '377.34'
$ws.Range('E5').Value = '  +6.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.79'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -2.22%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.97'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.41'
$ws.Range('E13').Value = '  -2.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.396.55'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.938.61'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.944'
$ws.Range('E17').Value = '  -5.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.500.15'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.02'
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0950'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.40'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '262.26'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.81'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.34'
$ws.Range('E26').Value = '  +17.00%  '
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  -5.47%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.35'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.84'
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('E32').Value = '  -6.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.84'
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.76'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '34.13'
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  -6.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.09'
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.61'
$ws.Range('E41').Value = '  -6.61%  '
$ws.Range('E42').Value = '  -5.31%  '
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '124.28'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.77'
$ws.Range('E45').Value = '  -5.12%  '
$ws.Range('E46').Value = '  -4.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.276'
$ws.Range('E47').Value = '  +14.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.021.75'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.18'
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('E51').Value = '  -2.31%  '
